# Update web content to display table views.
#
# The RS0004 description/performance block is renumbered: the redundant
# "RS0004" path segment is dropped from every ASHRAE205.RS_instance.* label,
# which shifts every subsequent data-group/data-element row up by one, and a
# new "-" units placeholder is added for cycling_degradation_coefficient.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Data-group path labels (column A)
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "ASHRAE205.RS_instance.description"
$ws.Range("A16").Value = "ASHRAE205.RS_instance.description.product_information"
$ws.Range("A17").ClearContents()

$ws.Range("A23").Value = "ASHRAE205.RS_instance.performance"
$ws.Range("A24").ClearContents()

$ws.Range("A26").Value = "ASHRAE205.RS_instance.performance.performance_map_cooling"
$ws.Range("A27").Value = "ASHRAE205.RS_instance.performance.performance_map_standby"
$ws.Range("A28").ClearContents()

# ---------------------------------------------------------------------
# 2) product_information data elements (rows 17-22), shifted up one row
#    and re-indented (16 -> 12 leading spaces)
# ---------------------------------------------------------------------
$ws.Range("B17").Value = "            outdoor_unit_manufacturer"
$ws.Range("C17").Value = "ColdAirInc"

$ws.Range("B18").Value = "            outdoor_unit_model_number"
$ws.Range("C18").Value = "AC.2019.01"

$ws.Range("B19").Value = "            indoor_unit_manufacturer"
$ws.Range("C19").Value = "ColdAirInc"

$ws.Range("B20").Value = "            indoor_unit_model_number"
$ws.Range("C20").Value = "AC.2019.02"

$ws.Range("B21").Value = "            refrigerant_type"
$ws.Range("C21").Value = "R410A"

$ws.Range("B22").Value = "            compressor_type"
$ws.Range("C22").Value = "SCROLL"

$ws.Range("B23").ClearContents()
$ws.Range("C23").ClearContents()

# ---------------------------------------------------------------------
# 3) performance data elements (rows 24-25), shifted up one row and
#    re-indented (12 -> 8 leading spaces)
# ---------------------------------------------------------------------
$ws.Range("B24").Value = "        compressor_control_method"
$ws.Range("C24").Value = "DYNAMIC"

$ws.Range("B25").Value = "        cycling_degradation_coefficient"
$ws.Range("C25").Value = 0.25
$ws.Range("D25").Value = "-"

$ws.Range("B26").ClearContents()
$ws.Range("C26").Value = "$" + "performance_map_cooling"

$ws.Range("C27").Value = "$" + "performance_map_standby"
$ws.Range("C28").ClearContents()

# ---------------------------------------------------------------------
# 4) "Required" column (E) - clear the checkmarks that no longer apply
#    to product_information / performance group members
# ---------------------------------------------------------------------
$ws.Range("E18").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("E24").ClearContents()
$ws.Range("E25").ClearContents()

# Drop the now-empty trailing row so the used range shrinks back down.
$ws.Rows.Item(28).Delete()

# ---------------------------------------------------------------------
# 5) Comments on sheet 1 - re-point / re-word to match the shifted rows
# ---------------------------------------------------------------------
$ws.Range("B18").Comment.Text("Model number of the outdoor unit")
$ws.Range("B19").Comment.Text("Name of the indoor unit manufacturer")
$ws.Range("B20").Comment.Text("Model number of the indoor unit")
$ws.Range("B21").Comment.Text("Type of refrigerant")
$ws.Range("B22").Comment.Text("Type of compressor")
$ws.Range("A27").Comment.Text("Data group describing standby performance")

$ws.Range("B23").Comment.Delete()
$ws.Range("B26").Comment.Delete()
$ws.Range("A28").Comment.Delete()

$ws.Range("A16").AddComment("Data group describing product information")
$ws.Range("B17").AddComment("Name of the outdoor unit manufacturer")
$ws.Range("B25").AddComment("Cycling degradation coefficient (CD) as described in AHRI 550/590 or AHRI 551/591")
$ws.Range("A26").AddComment("Data group describing cooling performance over a range of conditions")

# ---------------------------------------------------------------------
# 6) Data validations - RS_ID list gains RS0005/RS0006; the compressor
#    type / control method validations live one row higher now
# ---------------------------------------------------------------------
$ws.Range("C6").Validation.Formula1 = '"RS0001,RS0002,RS0003,RS0004,RS0005,RS0006"'

$ws.Range("C23").Validation.Delete()
$ws.Range("C22").Validation.Add(3, 1, 1, '"RECIPROCATING,SCREW,CENTRIFUGAL,ROTARY,SCROLL"')

$ws.Range("C25").Validation.Delete()
$ws.Range("C24").Validation.Add(3, 1, 1, '"STAGED,DYNAMIC"')

# ---------------------------------------------------------------------
# 7) performance_map_cooling / performance_map_standby sheets - drop the
#    redundant "RS0004" path segment from the title cell, and update the
#    standby sheet's "lookup variables" comment wording
# ---------------------------------------------------------------------
$wsCooling = $wb.Worksheets.Item("performance_map_cooling")
$wsCooling.Range("A1").Value = "ASHRAE205.RS_instance.performance.performance_map_cooling.grid_variables"

$wsStandby = $wb.Worksheets.Item("performance_map_standby")
$wsStandby.Range("A1").Value = "ASHRAE205.RS_instance.performance.performance_map_standby.grid_variables"
$wsStandby.Range("B2").Comment.Text("Data group defining the lookup variables for standby performance")
